$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 65
$ws.Range("C2").Value = 14
$ws.Range("D2").Value = 113
$ws.Range("E2").Value = 132
$ws.Range("I2").Value = 1
$ws.Range("K2").Value = 13
$ws.Range("L2").Value = 4
$ws.Range("M2").Value = 10
